$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Fat (%)"
$ws.Range("C1").Value = "Weight"

$ws.Range("H19").Select()
